# "updating ppt to other factors"
#
# 1) Insert a brand-new "Other factors:" slide (Title and Content layout)
#    right before the deck's final (design-element) slide, so it becomes
#    the new slide 11 and the old slide 11 shifts down to slide 12.
# 2) Re-cache the auto date field ("datetimeFigureOut") on the slide
#    master + every slide layout from 6/4/2022 -> 6/6/2022.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. New slide
# ---------------------------------------------------------------------
$layouts = $p.SlideMaster.CustomLayouts
$titleAndContent = $layouts.Item(2)   # "Title and Content" - same layout every other slide uses
$newSlide = $p.Slides.AddSlide(11, $titleAndContent)

# Title placeholder
$title = $newSlide.Shapes.Item(1)
$titleTr = $title.TextFrame.TextRange
$titleTr.Text = "Other factors:"
$titleTr.LanguageID = "en-IN"

# Body / content placeholder
$body = $newSlide.Shapes.Item(2)
$bodyTr = $body.TextFrame.TextRange
$bodyTr.Text = "State : huge variance in chargeoff% from state to state.`rZip code : huge variance in chargeoff% from one to other.`rHouse category: others have high chargeoff%`r`r`r`r`r`rChargeoff% = proportion of charged off members across total sampling."
$bodyTr.LanguageID = "en-IN"

# Shrink text on overflow (normAutofit)
$body.TextFrame.AutoSize = 2

# Last four paragraphs (two blank spacer lines + the "Chargeoff% = ..."
# footnote line) have their bullet removed.
for ($i = 6; $i -le 9; $i++) {
    $para = $bodyTr.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Type = 0
}

# ---------------------------------------------------------------------
# 2. Re-cache "datetimeFigureOut" field text: 6/4/2022 -> 6/6/2022
#    (slide master + all custom layouts)
# ---------------------------------------------------------------------
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lyt = $layouts.Item($i)
    for ($j = 1; $j -le $lyt.Shapes.Count; $j++) {
        $sh = $lyt.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "6/6/2022"
        }
    }
}

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "6/6/2022"
    }
}

Write-Output "Slide count: $($p.Slides.Count)"
